# The deck has a long run of slide-1 "time" textboxes (shape name "time")
# whose labels were re-formatted from the legacy "95 Dec 120505"-style
# stamp into ISO-8601 ("1995-12-12T05:05:00Z"). Walk every shape on the
# (only) slide and rewrite any textbox whose text matches the legacy
# pattern into the corresponding ISO-8601 timestamp.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$count = $s.Shapes.Count
for ($i = 1; $i -le $count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -eq $false) {
        continue
    }
    if ($shp.Name -ne "time") {
        continue
    }

    $tr = $shp.TextFrame.TextRange
    $old = $tr.Text
    if ($old -match '^95 Dec 12(\d{2})(\d{2})$') {
        $hh = $matches[1]
        $mm = $matches[2]
        $new = "1995-12-12T" + $hh + ":" + $mm + ":00Z"
        $tr.Text = $new
    }
}
